$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp / value data for rows 2 through 42 (columns A, B, C)
$data = @(
    @(45825,                 5.948,               10.463),
    @(45825.01041666666,     8.701000000000001,   0.254),
    @(45825.02083333334,     0.485,               7.58),
    @(45825.03125,           1.578,               3.483),
    @(45825.04166666666,     3.551,               0),
    @(45825.05208333334,     3.374,               0),
    @(45825.0625,            7.682,               0),
    @(45825.07291666666,     21.438,              0),
    @(45825.08333333334,     21.775,              0),
    @(45825.09375,           13.267,              0),
    @(45825.10416666666,     19.599,              0),
    @(45825.11458333334,     34.697,              0),
    @(45825.125,             34.327,              0),
    @(45825.13541666666,     20.304,              0),
    @(45825.14583333334,     35.943,              0),
    @(45825.15625,           9.301,               0),
    @(45825.16666666666,     1.788,               1.235),
    @(45825.17708333334,     0.708,               1.633),
    @(45825.1875,            0.013,               4.135),
    @(45825.19791666666,     7.906,               0.267),
    @(45825.20833333334,     5.255,               0),
    @(45825.21875,           3.173,               0),
    @(45825.22916666666,     7.908,               0.153),
    @(45825.23958333334,     8.369,               0.446),
    @(45825.25,              19.66,               0.019),
    @(45825.26041666666,     2.818,               4.206),
    @(45825.27083333334,     0,                   8.301),
    @(45825.28125,           0,                   10.119),
    @(45825.29166666666,     15.06,               1.412),
    @(45825.30208333334,     6.107,               2.326),
    @(45825.3125,            0,                   3.164),
    @(45825.32291666666,     0.155,               1.079),
    @(45825.33333333334,     20.104,              0),
    @(45825.34375,           0.367,               8.58),
    @(45825.35416666666,     0.008,               6.079),
    @(45825.36458333334,     0,                   13.236),
    @(45825.375,             0,                   8.914),
    @(45825.38541666666,     0,                   10.834),
    @(45825.39583333334,     0,                   3.232),
    @(45825.40625,           0,                   0.916),
    @(45825.41666666666,     3.231,               3.497)
)

$startRow = 2

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row[0]
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
